$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 112, shifting existing rows 112-174 down to 113-175
# (mirrors the "dimension A1:R174 -> A1:R175" + shifted data seen in the diff).
$ws.Rows.Item(112).Insert()

# Populate the newly inserted row 112 with the new weekly price record.
$ws.Range("A112").Value = 9
$ws.Range("B112").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C112").Value = "Metropolitana"
$ws.Range("D112").Value = 45126
$ws.Range("E112").Value = 13
$ws.Range("F112").Value = 100112022
$ws.Range("G112").Value = "Arveja Verde"
$ws.Range("H112").Value = "Perfection"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 52
$ws.Range("K112").Value = 30000
$ws.Range("L112").Value = 32000
$ws.Range("M112").Value = 31000
$ws.Range("N112").Value = "`$/malla 25 kilos"
$ws.Range("O112").Value = "Provincia de Huasco"
$ws.Range("P112").Value = 1240
$ws.Range("Q112").Value = 25
$ws.Range("R112").Value = "Hortaliza"
